$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.688.99'
$ws.Range("D3").Value = '3.937.20'
$ws.Range("E3").Value = '  +0.72%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '531.52'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +8.99%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.84'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.09%  '
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.732'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.66%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.174'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.39%  '
$ws.Range("E11").Value = '  -0.28%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '43.00'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.22%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.50'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.95%  '
$ws.Range("D14").Value = '4.567.80'
$ws.Range("E14").Value = '  +0.84%  '
$ws.Range("D15").Value = '3.955.21'
$ws.Range("E15").Value = '  +1.19%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.11'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.27%  '
$ws.Range("E17").Value = '  -0.10%  '
$ws.Range("E18").Value = '  +7.15%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '19.92'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.55%  '
$ws.Range("D20").Value = '69.612.11'
$ws.Range("E20").Value = '  +1.96%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '435.98'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.25%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.42'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.35%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.60'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.19%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '88.50'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.26%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.07'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +12.61%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.94'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.26%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.02'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.56%  '
$ws.Range("E28").Value = '  -3.62%  '
$ws.Range("E29").Value = '  -1.29%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '705.56'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.85%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '13.41'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.36%  '
$ws.Range("E32").Value = '  -1.76%  '
$ws.Range("E33").Value = '  -1.43%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '68.55'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +13.65%  '
$ws.Range("E35").Value = '  +8.44%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.12'
$ws.Range("D36").Style = "Normal"
$ws.Range("E37").Value = '  -0.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '40.61'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.92%  '
$ws.Range("E39").Value = '  +0.62%  '
$ws.Range("E40").Value = '  +0.03%  '
$ws.Range("E41").Value = '  -0.06%  '
$ws.Range("E42").Value = '  +1.35%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.83'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.57%  '
$ws.Range("E44").Value = '  +6.43%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.02'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.78%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.20'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +14.13%  '
$ws.Range("E47").Value = '  +2.42%  '
$ws.Range("E48").Value = '  +1.62%  '
$ws.Range("D49").Value = '0.0₆0365'
$ws.Range("E49").Value = '  +6.54%  '
$ws.Range("E50").Value = '  -1.29%  '
$ws.Range("E51").Value = '  -0.97%  '
